$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out old content
$ws.Cells.Clear()

# Set new content
$ws.Range("A7").Value = "dooosraa"

# Adjust column A width
$ws.Columns.Item(1).ColumnWidth = 66.1796875

# Select A7 as active cell
$ws.Range("A7").Select()
